$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C11").Value = "Los Lagos"
$ws.Range("D11").Value = 44488
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 300000000
$ws.Range("G11").Value = "Espárragos"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 600
$ws.Range("K11").Value = 1700
$ws.Range("L11").Value = 1800
$ws.Range("M11").Value = 1750
$ws.Range("N11").Value = "$/kilo"
$ws.Range("O11").Value = "Provincia de Linares"
$ws.Range("P11").Value = 1750
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = "Hortaliza"
